$d = $word.ActiveDocument

# --- Edit 1: mark the (previously un-marked) inline picture run as
# "do not spell-check" (<w:noProof/>) -- the picture that sits right
# before the "...by creating a compl..." paragraph's lastRenderedPageBreak.
$shp = $d.InlineShapes(2)
$shp.Range.NoProofing = $true

# --- Edit 2: remove one of the three consecutive empty paragraphs that
# sit right before the "Quality assurance criteria /2" heading.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*Quality assurance criteria /2*") {
        $target = $i
        break
    }
}
if ($target -ne $null) {
    $empty = $d.Paragraphs($target - 1)
    $empty.Range.Delete()
}

# --- Edit 3: add a new (unused) numbered-list definition to numbering.xml
# -- mirrors a list style that was applied and then removed again, which
# leaves the abstractNum/num entries behind (numId=3 / abstractNumId=2).
$r = $d.Content
$r.Collapse(0)
$r.InsertParagraphAfter()
$lastIndex = $d.Paragraphs.Count
$scratch = $d.Paragraphs($lastIndex)
$scratch.Range.ListFormat.ApplyNumberDefault()
$scratch.Range.Delete()

Write-Output "done"
